# "kind of working speed limit sign"
# Update the simulation data: refresh the Statistics sheet with a new
# batch of sensor readings (9 new rows appended, first 6 rows' values
# replaced) and clear out the (now stale / unconfirmed) Accidents entry,
# leaving just the header row behind.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Statistics sheet
# ---------------------------------------------------------------------
$stats = $wb.Worksheets.Item("Statistics")

$rows = @(
    @("2024-08-26 15:04:47", 44.13390296570707, 5),
    @("2024-08-26 15:04:49", 42.16621551850734, 8),
    @("2024-08-26 15:04:51", 42.38298800889578, 9),
    @("2024-08-26 15:04:53", 44.16098196112721, 11),
    @("2024-08-26 15:06:01", 42.22672922880924, 12),
    @("2024-08-26 15:06:12", 42.15700632867082, 12),
    @("2024-08-26 15:06:17", 42.77970146747748, 13),
    @("2024-08-26 15:06:19", 42.73019503803481, 13),
    @("2024-08-26 15:09:01", 43.33023722596015, 14),
    @("2024-08-26 15:09:03", 43.39401072519414, 14),
    @("2024-08-26 15:09:05", 43.31036399285087, 14),
    @("2024-08-26 15:09:40", 43.11456799213904, 14),
    @("2024-08-26 15:09:42", 42.98245919569164, 14),
    @("2024-08-26 15:09:44", 38.53016086544418, 14),
    @("2024-08-26 15:09:46", 42.97772545656818, 17)
)

$r = 2
foreach ($row in $rows) {
    $stats.Cells.Item($r, 1).Value = $row[0]
    $stats.Cells.Item($r, 2).Value = $row[1]
    $stats.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# ---------------------------------------------------------------------
# Accidents sheet - remove the single recorded accident row, leaving
# only the header row in place.
# ---------------------------------------------------------------------
$accidents = $wb.Worksheets.Item("Accidents")
$accidents.Rows.Item(2).Delete()
